# The "municipio-nombre" column (L) was re-classified from a "measure" to a
# (curated) "dimension", matching the pattern already used by the
# "provincia-nombre" (M) and "comarca-nombre" (N) columns.
#
#   L2: iaest-measure:municipio-nombre  ->  sdmx-dimension:refArea
#   L3: medida                          ->  dim
#   L4: xsd:int                         ->  URI-Municipio

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "sdmx-dimension:refArea"
$ws.Range("L3").Value = "dim"
$ws.Range("L4").Value = "URI-Municipio"
